# chore: update Sheets via scheduled runner
# Refreshes cached market-price derived columns (H:N) for a handful of
# Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1012208.75
$ws.Range("I88").Value = 2224959.5
$ws.Range("J88").Value = 1583.1666
$ws.Range("K88").Value = 2224959.5
$ws.Range("L88").Value = 1583.1666
$ws.Range("M88").Value = -2224553.5
$ws.Range("N88").Value = -2395.1666
$ws.Range("H91").Value = 1012208.75
$ws.Range("I91").Value = 2224959.5
$ws.Range("J91").Value = 1583.1666
$ws.Range("K91").Value = 2224959.5
$ws.Range("L91").Value = 1583.1666
$ws.Range("M91").Value = -2223555.5
$ws.Range("N91").Value = -4391.1666
$ws.Range("H113").Value = 244216.72
$ws.Range("I113").Value = 385881.72
$ws.Range("J113").Value = 2552.8823
$ws.Range("K113").Value = 385881.72
$ws.Range("L113").Value = 2552.8823
$ws.Range("M113").Value = -382627.72
$ws.Range("N113").Value = -9060.882300000001
$ws.Range("H116").Value = 2310053.8
$ws.Range("I116").Value = 11906978
$ws.Range("J116").Value = 6791.96
$ws.Range("K116").Value = 11906978
$ws.Range("L116").Value = 6791.96
$ws.Range("M116").Value = -11903536
$ws.Range("N116").Value = -13675.96
$ws.Range("H129").Value = 1429.0975
$ws.Range("I129").Value = 652.6667
$ws.Range("K129").Value = 1958.0001
$ws.Range("M129").Value = 3041.9999
$ws.Range("H132").Value = 1931813
$ws.Range("I132").Value = 2101937.8
$ws.Range("J132").Value = 3731.6667
$ws.Range("K132").Value = 6305813.399999999
$ws.Range("L132").Value = 11195.0001
$ws.Range("M132").Value = -6303283.399999999
$ws.Range("N132").Value = -16255.0001
$ws.Range("H135").Value = 2585.7778
$ws.Range("I135").Value = 2693.875
$ws.Range("J135").Value = 2276.9285
$ws.Range("K135").Value = 24244.875
$ws.Range("L135").Value = 20492.3565
$ws.Range("M135").Value = -21709.875
$ws.Range("N135").Value = -25562.3565
$ws.Range("H137").Value = 1074.2609
$ws.Range("I137").Value = 844.08105
$ws.Range("J137").Value = 2020.5555
$ws.Range("K137").Value = 2532.24315
$ws.Range("L137").Value = 6061.666499999999
$ws.Range("M137").Value = 17.75684999999976
$ws.Range("N137").Value = -11161.6665

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5694.5845
$ws.Range("I32").Value = 2869.5774
$ws.Range("K32").Value = 2869.5774
$ws.Range("M32").Value = -2582.5774
$ws.Range("H61").Value = 971.6977000000001
$ws.Range("I61").Value = 886.32434
$ws.Range("J61").Value = 1498.1666
$ws.Range("K61").Value = 886.32434
$ws.Range("L61").Value = 1498.1666
$ws.Range("M61").Value = -674.32434
$ws.Range("N61").Value = -1922.1666
$ws.Range("H110").Value = 1346.4445
$ws.Range("I110").Value = 841.4737
$ws.Range("J110").Value = 2545.75
$ws.Range("K110").Value = 841.4737
$ws.Range("L110").Value = 2545.75
$ws.Range("M110").Value = 1203.5263
$ws.Range("N110").Value = -6635.75
$ws.Range("H132").Value = 974.42426
$ws.Range("I132").Value = 779.94446
$ws.Range("J132").Value = 1849.5834
$ws.Range("K132").Value = 2339.83338
$ws.Range("L132").Value = 5548.7502
$ws.Range("M132").Value = 190.16662
$ws.Range("N132").Value = -10608.7502
$ws.Range("H136").Value = 971.6977000000001
$ws.Range("I136").Value = 886.32434
$ws.Range("J136").Value = 1498.1666
$ws.Range("K136").Value = 2658.97302
$ws.Range("L136").Value = 4494.4998
$ws.Range("M136").Value = -108.9730199999999
$ws.Range("N136").Value = -9594.4998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1640
$ws.Range("I86").Value = 1437.5
$ws.Range("J86").Value = 2450
$ws.Range("K86").Value = 1437.5
$ws.Range("L86").Value = 2450
$ws.Range("M86").Value = -314.5
$ws.Range("N86").Value = -4696
$ws.Range("H89").Value = 1640
$ws.Range("I89").Value = 1437.5
$ws.Range("J89").Value = 2450
$ws.Range("K89").Value = 7187.5
$ws.Range("L89").Value = 12250
$ws.Range("M89").Value = -1571.5
$ws.Range("N89").Value = -23482
$ws.Range("H130").Value = 40624
$ws.Range("J130").Value = 40624
$ws.Range("L130").Value = 40624
$ws.Range("N130").Value = -50664
$ws.Range("H134").Value = 1856.8182
$ws.Range("I134").Value = 1566.8948
$ws.Range("J134").Value = 2250.2856
$ws.Range("K134").Value = 4700.6844
$ws.Range("L134").Value = 6750.8568
$ws.Range("M134").Value = -2165.6844
$ws.Range("N134").Value = -11820.8568

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5733.3335
$ws.Range("I16").Value = 3033.3333
$ws.Range("K16").Value = 3033.3333
$ws.Range("M16").Value = -2746.3333
$ws.Range("H22").Value = 432
$ws.Range("I22").Value = 335
$ws.Range("K22").Value = 335
$ws.Range("M22").Value = 15
$ws.Range("H31").Value = 1445.0625
$ws.Range("I31").Value = 989.6857
$ws.Range("J31").Value = 4632.7
$ws.Range("K31").Value = 989.6857
$ws.Range("L31").Value = 4632.7
$ws.Range("M31").Value = -694.6857
$ws.Range("N31").Value = -5222.7
$ws.Range("H34").Value = 1445.0625
$ws.Range("I34").Value = 989.6857
$ws.Range("J34").Value = 4632.7
$ws.Range("K34").Value = 989.6857
$ws.Range("L34").Value = 4632.7
$ws.Range("M34").Value = -787.6857
$ws.Range("N34").Value = -5036.7
$ws.Range("H99").Value = 2805
$ws.Range("I99").Value = 2375.3845
$ws.Range("J99").Value = 4666.6665
$ws.Range("K99").Value = 2375.3845
$ws.Range("L99").Value = 4666.6665
$ws.Range("M99").Value = -877.3845000000001
$ws.Range("N99").Value = -7662.6665
$ws.Range("H113").Value = 5733.3335
$ws.Range("I113").Value = 3033.3333
$ws.Range("K113").Value = 3033.3333
$ws.Range("M113").Value = -863.3332999999998
$ws.Range("H126").Value = 2805
$ws.Range("I126").Value = 2375.3845
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 7126.1535
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -4656.1535
$ws.Range("N126").Value = -18939.9995
$ws.Range("H132").Value = 1604.3939
$ws.Range("I132").Value = 1220.2593
$ws.Range("J132").Value = 3333
$ws.Range("K132").Value = 3660.7779
$ws.Range("L132").Value = 9999
$ws.Range("M132").Value = -1130.7779
$ws.Range("N132").Value = -15059
$ws.Range("H134").Value = 1316.5358
$ws.Range("I134").Value = 1093.9524
$ws.Range("J134").Value = 1984.2858
$ws.Range("K134").Value = 3281.857199999999
$ws.Range("L134").Value = 5952.857400000001
$ws.Range("M134").Value = -746.8571999999995
$ws.Range("N134").Value = -11022.8574

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1190847.9
$ws.Range("I122").Value = 294.94595
$ws.Range("J122").Value = 10000939
$ws.Range("K122").Value = 2654.51355
$ws.Range("L122").Value = 90008451
$ws.Range("M122").Value = -204.5135499999997
$ws.Range("N122").Value = -90013351

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2598.5356
$ws.Range("I102").Value = 1381.0555
$ws.Range("J102").Value = 4790
$ws.Range("K102").Value = 1381.0555
$ws.Range("L102").Value = 4790
$ws.Range("M102").Value = 240.9445000000001
$ws.Range("N102").Value = -8034
$ws.Range("H122").Value = 2240.4
$ws.Range("I122").Value = 2057.7144
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 6173.1432
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -3723.1432
$ws.Range("N122").Value = -12900.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 318.47827
$ws.Range("I22").Value = 311.25
$ws.Range("J22").Value = 366.66666
$ws.Range("K22").Value = 311.25
$ws.Range("L22").Value = 366.66666
$ws.Range("M22").Value = -16.25
$ws.Range("N22").Value = -956.66666
$ws.Range("H27").Value = 318.47827
$ws.Range("I27").Value = 311.25
$ws.Range("J27").Value = 366.66666
$ws.Range("K27").Value = 311.25
$ws.Range("L27").Value = 366.66666
$ws.Range("M27").Value = -204.25
$ws.Range("N27").Value = -580.66666
$ws.Range("H136").Value = 1437.9231
$ws.Range("I136").Value = 686.8570999999999
$ws.Range("J136").Value = 2314.1667
$ws.Range("K136").Value = 2060.5713
$ws.Range("L136").Value = 6942.500100000001
$ws.Range("M136").Value = 489.4287000000004
$ws.Range("N136").Value = -12042.5001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 559355.5600000001
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 559355.5600000001
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 559355.5600000001
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -559579.5600000001
$ws.Range("H113").Value = 435.2143
$ws.Range("J113").Value = 473.375
$ws.Range("L113").Value = 1420.125
$ws.Range("N113").Value = -5760.125
$ws.Range("H132").Value = 1331.3243
$ws.Range("I132").Value = 1309.3334
$ws.Range("J132").Value = 1360.1875
$ws.Range("K132").Value = 3928.0002
$ws.Range("L132").Value = 4080.5625
$ws.Range("M132").Value = -1398.0002
$ws.Range("N132").Value = -9140.5625
$ws.Range("H136").Value = 994.8679
$ws.Range("I136").Value = 827.5161000000001
$ws.Range("J136").Value = 1230.6818
$ws.Range("K136").Value = 2482.5483
$ws.Range("L136").Value = 3692.0454
$ws.Range("M136").Value = 67.45169999999962
$ws.Range("N136").Value = -8792.045399999999
